$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 15.0677871054159
$ws.Range("E2").Value = 0.8086528589963745
$ws.Range("D3").Value = 13.65273615994774
$ws.Range("E3").Value = 2.056412934489027
$ws.Range("D4").Value = 15.84884973754504
$ws.Range("E4").Value = 0.4737798588840921
$ws.Range("D5").Value = 14.82296216149963
$ws.Range("E5").Value = 1.470271652439198
$ws.Range("D6").Value = 16.27049478550594
$ws.Range("E6").Value = 0.5175175310755836
$ws.Range("D7").Value = 15.36140106242473
$ws.Range("E7").Value = 1.160568058301429
$ws.Range("D8").Value = 16.27619821733979
$ws.Range("E8").Value = 0.4962008151841023
$ws.Range("D9").Value = 15.48154308277698
$ws.Range("E9").Value = 0.7801327716119801
$ws.Range("D10").Value = 16.82246103450655
$ws.Range("E10").Value = 0.521185082417052
$ws.Range("D11").Value = 16.13016089284572
$ws.Range("E11").Value = 0.8896458082423097
$ws.Range("D12").Value = 20.04140488260402
$ws.Range("E12").Value = 2.00159589401732
$ws.Range("D13").Value = 17.34250682055415
$ws.Range("E13").Value = 0.5344189831229733
$ws.Range("D14").Value = 20.93593117027743
$ws.Range("E14").Value = 2.659476564274648
$ws.Range("D15").Value = 17.63861984181357
$ws.Range("E15").Value = 1.00072081656106
$ws.Range("D16").Value = 21.80860802039622
$ws.Range("E16").Value = 2.621958038061083
$ws.Range("D17").Value = 17.41312086157644
$ws.Range("E17").Value = 1.098932289689143
$ws.Range("D18").Value = 22.07572639411203
$ws.Range("E18").Value = 2.706390963726625
$ws.Range("D19").Value = 18.18954593210566
$ws.Range("E19").Value = 1.014170911387849
